# Updates the "Estado de Cuenta" worksheet:
#  - Replaces the worker JAWIN MERCADO PADILLA (1052740323) with
#    LUIS ALFONSO MORALES PAREDES (1052982551) on the first data row,
#    and renumbers the "Periodo Mora" column sequentially (2410..2508)
#    for the now-single worker, adding a new period 2508 (part 1 of
#    new estado de cuenta).
#  - Updates the totals (Valor Mora, Cant. Trabajadores, Cant. Periodos).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Totals block
$ws.Range("E11").Value = 525200
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 11

# Data rows 16-26: single worker, sequential periods 2410..2508
$ws.Range("C16").Value = "1052982551"
$ws.Range("D16").Value = "LUIS ALFONSO MORALES PAREDES"
$ws.Range("E16").Value = "2410"
$ws.Range("F16").Value = 5200
$ws.Range("G16").Value = 1300000

$ws.Range("C17").Value = "1052982551"
$ws.Range("D17").Value = "LUIS ALFONSO MORALES PAREDES"
$ws.Range("E17").Value = "2411"
$ws.Range("F17").Value = 52000
$ws.Range("G17").Value = 1300000

$ws.Range("C18").Value = "1052982551"
$ws.Range("D18").Value = "LUIS ALFONSO MORALES PAREDES"
$ws.Range("E18").Value = "2412"
$ws.Range("F18").Value = 52000
$ws.Range("G18").Value = 1300000

$ws.Range("C19").Value = "1052982551"
$ws.Range("D19").Value = "LUIS ALFONSO MORALES PAREDES"
$ws.Range("E19").Value = "2501"
$ws.Range("F19").Value = 52000
$ws.Range("G19").Value = 1300000

$ws.Range("C20").Value = "1052982551"
$ws.Range("D20").Value = "LUIS ALFONSO MORALES PAREDES"
$ws.Range("E20").Value = "2502"
$ws.Range("F20").Value = 52000
$ws.Range("G20").Value = 1300000

$ws.Range("C21").Value = "1052982551"
$ws.Range("D21").Value = "LUIS ALFONSO MORALES PAREDES"
$ws.Range("E21").Value = "2503"
$ws.Range("F21").Value = 52000
$ws.Range("G21").Value = 1300000

$ws.Range("C22").Value = "1052982551"
$ws.Range("D22").Value = "LUIS ALFONSO MORALES PAREDES"
$ws.Range("E22").Value = "2504"
$ws.Range("F22").Value = 52000
$ws.Range("G22").Value = 1300000

$ws.Range("C23").Value = "1052982551"
$ws.Range("D23").Value = "LUIS ALFONSO MORALES PAREDES"
$ws.Range("E23").Value = "2505"
$ws.Range("F23").Value = 52000
$ws.Range("G23").Value = 1300000

$ws.Range("C24").Value = "1052982551"
$ws.Range("D24").Value = "LUIS ALFONSO MORALES PAREDES"
$ws.Range("E24").Value = "2506"
$ws.Range("F24").Value = 52000
$ws.Range("G24").Value = 1300000

$ws.Range("C25").Value = "1052982551"
$ws.Range("D25").Value = "LUIS ALFONSO MORALES PAREDES"
$ws.Range("E25").Value = "2507"
$ws.Range("F25").Value = 52000
$ws.Range("G25").Value = 1300000

$ws.Range("C26").Value = "1052982551"
$ws.Range("D26").Value = "LUIS ALFONSO MORALES PAREDES"
$ws.Range("E26").Value = "2508"
$ws.Range("F26").Value = 52000
$ws.Range("G26").Value = 1300000

# Columns are best-fit; re-autofit now that the longer name occupies every row
$ws.Columns("B:J").AutoFit() | Out-Null
